$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 13 (pushes old row14 "Success" row down to row15)
$ws.Rows("13").Insert()

# --- Row 2 ---
$ws.Range("C2").Value = 13
$ws.Range("H2").Value = 0

# --- Row 3 ---
$ws.Range("C3").Value = 1

# --- Row 5 ---
$ws.Range("C5").Value = 1

# --- Row 6 ---
$ws.Range("C6").Value = 1

# --- Row 8: formula replaced with a plain value ---
$ws.Range("C8").Value = -210

# --- Row 9: Int/Wis pattern ---
$ws.Range("B9").Value2 = $ws.Range("B2").Value2
$ws.Range("D9").Value = 20
$ws.Range("C9").Formula = "=(C2)*D9"

# --- Row 10: Spell level pattern ---
$ws.Range("B10").Value2 = $ws.Range("B3").Value2
$ws.Range("D10").Value = -20
$ws.Range("C10").Formula = "=D10*(C3)"

# --- Row 11: Penalty pattern ---
$ws.Range("B11").Value2 = $ws.Range("B4").Value2
$ws.Range("D11").Value = -5
$ws.Range("C11").Formula = "=D11*C4"

# --- Row 12: Unit level pattern ---
$ws.Range("B12").Value2 = $ws.Range("B5").Value2
$ws.Range("D12").Value = 20
$ws.Range("C12").Formula = "=D12*(C5)"

# --- Row 13 (new row): Skill pattern ---
$ws.Range("B13").Value2 = $ws.Range("B6").Value2
$ws.Range("D13").Value = 50
$ws.Range("C13").Formula = "=D13*C6"

# --- Row 15 (shifted from old row14): update SUM range to include new row14 ---
$ws.Range("C15").Formula = "=SUM(C8:C14)"

# --- Selection update ---
$ws.Range("C2").Select()

$wb.Save()
